# Apply "ajuste de admitidos e desligados" changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1)
$ws.Range("B1").Value = "saldo_ajuste"
$ws.Range("C1").Value = "admitidos_ajuste"
$ws.Range("D1").Value = "desligados_ajuste"

# New "admitidos_ajuste" (column C) values for rows 2-49
$cValues = @(199,187,179,96,96,116,195,187,232,273,200,145,250,267,266,207,242,261,226,296,299,272,244,198,308,273,266,311,266,277,204,392,297,245,401,189,303,265,409,292,395,343,295,318,374,369,299,210)

# New "desligados_ajuste" (column D) values for rows 2-49
$dValues = @(199,184,246,207,204,156,136,118,136,165,181,163,196,242,274,208,242,229,198,237,223,257,204,229,235,203,221,202,242,229,209,264,198,235,200,208,243,232,264,240,256,303,274,385,237,273,238,250)

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}

# Row 44 (competenciamov 202307) also has an adjusted "saldo_ajuste" value
$ws.Range("B44").Value = 21
